$d = $word.ActiveDocument

# 1. Update the "time out server-side" bullet to mention cancelling the
#    request client-side as well.
$old1 = "Models occasionally may time out server-side and the user won’t get notified. This could be removed by showing a message after a certain amount of time has passed."
$new1 = "Models occasionally may time out server-side and the user won’t get notified. This could be removed by showing a message after a certain amount of time has passed and having the client cancel the request."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# 2. Append four new bulleted list items after the last paragraph
#    ("Given certain inputs, ..."), reusing that paragraph's numbered-list
#    formatting.
$bullets = @(
    "When there is an excess of saved articles, the “more” button does not work. This can be resolved by removing the article limit and adding a scrollable container in its place",
    "The summary section is able to scroll, but there is no visible scrollbar. This can be fixed by having the container inherit the global scrollbar css",
    "Upload button does not work when editing the first summary. This can be fixed by turning off edit mode after uploading.",
    "“Article Submitted” text does not automatically dismiss itself. This can be fixed by adding its own timeouts."
)

foreach ($bulletText in $bullets) {
    $lastPara = $d.Paragraphs.Last
    $lastPara.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Last
    $newPara.Range.Text = $bulletText
}
